$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17/18: Litecoin and ShibaInu swap positions with updated values
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.35"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001074"
$ws.Range("E18").Value = "  +1.76%  "

# Update Price (D) and Volume(1h) (E) values for all other changed rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.440.22"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.40"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.76"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5324"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4053"
$ws.Range("E8").Value = "  +7.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07603"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.85"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.366"
$ws.Range("E12").Value = "  +4.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.95"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.559"
$ws.Range("E15").Value = "  +4.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.826.94"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06606"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.62"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.075"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.459.58"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.158"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.485"
$ws.Range("E26").Value = "  +8.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.03"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.57"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.035.60"
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.72"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.126"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1094"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.687"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.659"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07190"
$ws.Range("E35").Value = "  +12.36%  "
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02343"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.230"
$ws.Range("E38").Value = "  +4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.773"
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6279"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.185"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.398"
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.709"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5851"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.24"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.199"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  +0.83%  "
